# The presentation used to start with a "flow diagram" slide (Preprocess /
# Process / PostProcess / Conditions boxes) followed by two image slides.
# That diagram slide is now obsolete and is removed, so the two image
# slides shift up to become slide 1 and slide 2.

$p = $ppt.ActivePresentation

# Slide 1 is the outdated diagram slide - delete it. PowerPoint
# automatically renumbers/shifts the remaining slides (old slide 2 becomes
# the new slide 1, old slide 3 becomes the new slide 2) and updates the
# slide id list accordingly.
$p.Slides.Item(1).Delete()
